$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matches source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.451.41"
$ws.Range("E2").Value = "  -2.76%  "
$ws.Range("D3").Value = "2.248.88"
$ws.Range("E3").Value = "  -3.80%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").Value = "233.77"
$ws.Range("E5").Value = "  -2.20%  "
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("D7").Value = "69.67"
$ws.Range("E7").Value = "  -2.97%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -4.14%  "
$ws.Range("D10").Value = "0.0992"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").Value = "58.40"
$ws.Range("E11").Value = "  +0.61%  "
$ws.Range("D12").Value = "35.86"
$ws.Range("E12").Value = "  +11.25%  "
$ws.Range("E13").Value = "  -1.60%  "
$ws.Range("D14").Value = "6.77"
$ws.Range("E14").Value = "  -4.77%  "
$ws.Range("D15").Value = "2.581.17"
$ws.Range("E15").Value = "  -3.97%  "
$ws.Range("E16").Value = "  -6.65%  "
$ws.Range("E17").Value = "  -3.74%  "
$ws.Range("D18").Value = "2.245.44"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("D19").Value = "42.202.23"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").Value = "0.0₃0976"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").Value = "73.46"
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("D23").Value = "236.35"
$ws.Range("E23").Value = "  -6.05%  "
$ws.Range("E24").Value = "  +4.64%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "3.67"
$ws.Range("E26").Value = "  -1.80%  "
$ws.Range("D27").Value = "2.38"
$ws.Range("E27").Value = "  -4.41%  "
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("E29").Value = "  -3.89%  "
$ws.Range("D30").Value = "169.02"
$ws.Range("E30").Value = "  -3.67%  "
$ws.Range("E31").Value = "  -6.93%  "
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("E33").Value = "  -5.64%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.0719"
$ws.Range("E35").Value = "  -1.60%  "
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("E37").Value = "  -2.50%  "
$ws.Range("D38").Value = "21.59"
$ws.Range("E38").Value = "  +15.32%  "
$ws.Range("D39").Value = "2.27"
$ws.Range("E39").Value = "  -4.09%  "
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "5.96"
$ws.Range("E40").Value = "  -6.08%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "0.0268"
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").Value = "66.36"
$ws.Range("E42").Value = "  +2.41%  "
$ws.Range("D43").Value = "4.92"
$ws.Range("E43").Value = "  -6.07%  "
$ws.Range("D44").Value = "8.97"
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("E45").Value = "  -4.26%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.189"
$ws.Range("E46").Value = "  -3.15%  "
$ws.Range("B47").Value = "BinanceUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("B48").Value = "BitTorrent-New"
$ws.Range("C48").Value = "https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt"
$ws.Range("D48").Value = "0.0₃0156"
$ws.Range("E48").Value = "  +27.13%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "4.45"
$ws.Range("E49").Value = "  +11.03%  "
$ws.Range("E50").Value = "  -2.94%  "
$ws.Range("D51").Value = "2.35"
$ws.Range("E51").Value = "  -3.45%  "
